$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add the new "bubble diff" worksheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "bubble diff"

# --- Populate "bubble diff" in the same order the original author typed it  ---
# (this keeps the shared-string table index order identical to the target file)
$ws2.Range("D8").Value = "r"
$ws2.Range("F4").Value = "cm^5"
$ws2.Range("F5").Value = "um^5"
$ws2.Range("D4").Value = "Db"
$ws2.Range("D6").Value = "Qb"
$ws2.Range("F6").Value = "J/mol"
$ws2.Range("D2").Value = "T"
$ws2.Range("D1").Value = "R"

# --- Numeric / formula content (write all values & formulas before touching
#     number formats, so number-format application doesn't bleed its style
#     onto not-yet-created neighboring cells in the same row) ---
$ws2.Range("E1").Value = 8.3144600000000004
$ws2.Range("E2").Value = 3000

$ws2.Range("E4").Value = 1.19 * [Math]::Pow(10, -22)
$ws2.Range("E5").Formula = "=E4*100000000000000000000"
$ws2.Range("E6").Value = 418400

$ws2.Range("D10").Value = 0.002
$ws2.Range("E10").Formula = "=`$E`$5/D10^3*EXP(-`$E`$6/`$E`$1/`$E`$2)"

$ws2.Range("D11").Value = 0.003
$ws2.Range("E11").Formula = "=`$E`$5/D11^3*EXP(-`$E`$6/`$E`$1/`$E`$2)"

# --- Number formats (applied after all values/formulas are in place) ---
$ws2.Range("E4").NumberFormat = "0.00E+00"
$ws2.Range("E5").NumberFormat = "0.00E+00"
$ws2.Range("D9").NumberFormat = "0.00E+00"
$ws2.Range("D10").NumberFormat = "0.00E+00"
$ws2.Range("D11").NumberFormat = "0.00E+00"

# Column E is a touch wider than default, matching the source file's bestFit width
$ws2.Columns.Item(5).ColumnWidth = 11.33

# --- Selection bookkeeping: leave "bubble diff" remembering G30 ... ---
$ws2.Activate()
$ws2.Range("G30").Select()

# ...but return to Sheet1 as the active sheet/selection, matching the target
$ws1.Activate()
$ws1.Range("C19").Select()
